# global_variable_template.xlsx edit
# Fix global_variable_generator.py to import arrays properly for IO modules
# Added warning system
#
# This script reproduces, via Excel COM automation, the row insertions /
# cell edits described by the authoritative xml diff:
#   - Constants sheet: new alarm/VFD status bit, new CONST_* range-check
#     constants for air + water, and a new global alarm-clear flag.
#   - IO Mapping sheet: new MB_pumpN_faultRecord array mappings for each
#     of the three pumps (fault history word arrays).
#   - Minor view/selection + column-width cosmetics.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Constants" worksheet (sheet1.xml)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Constants")

# Work from the bottom of the sheet upward so that previously-computed
# target row numbers remain valid after each insert (an Insert() at a
# lower row never disturbs rows above it).

# -- Insert 8 new rows at 73..80 (CONST_* range-checking constants),
#    right after "ventilationCheckCounter" (row 72) and before the old
#    "tmpReal" row (old row 71, which becomes row 81).
$ws.Rows.Item("73:80").Insert()

$ws.Cells.Item(73, 2).Value = "CONST_AIR_TEMPERATURE_MAX"
$ws.Cells.Item(73, 3).Value = 1
$ws.Cells.Item(73, 4).Value = "WORD"
$ws.Cells.Item(73, 4).HorizontalAlignment = -4131
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(73, 6).Value = "x"
$ws.Cells.Item(73, 7).Value = "For range checking"

$ws.Cells.Item(74, 2).Value = "CONST_AIR_TEMPERATURE_MIN"
$ws.Cells.Item(74, 3).Value = 1
$ws.Cells.Item(74, 4).Value = "WORD"
$ws.Cells.Item(74, 4).HorizontalAlignment = -4131
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(74, 6).Value = "x"
$ws.Cells.Item(74, 7).Value = "For range checking"

$ws.Cells.Item(75, 2).Value = "CONST_AIR_HUMIDITY_MAX"
$ws.Cells.Item(75, 3).Value = 1
$ws.Cells.Item(75, 4).Value = "WORD"
$ws.Cells.Item(75, 4).HorizontalAlignment = -4131
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(75, 6).Value = "x"
$ws.Cells.Item(75, 7).Value = "For range checking"

$ws.Cells.Item(76, 2).Value = "CONST_AIR_HUMIDITY_MIN"
$ws.Cells.Item(76, 3).Value = 1
$ws.Cells.Item(76, 4).Value = "WORD"
$ws.Cells.Item(76, 4).HorizontalAlignment = -4131
$ws.Cells.Item(76, 5).Value = 0
$ws.Cells.Item(76, 6).Value = "x"
$ws.Cells.Item(76, 7).Value = "For range checking"

$ws.Cells.Item(77, 2).Value = "CONST_AIR_CO2_MAX"
$ws.Cells.Item(77, 3).Value = 1
$ws.Cells.Item(77, 4).Value = "WORD"
$ws.Cells.Item(77, 4).HorizontalAlignment = -4131
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(77, 6).Value = "x"
$ws.Cells.Item(77, 7).Value = "For range checking"

$ws.Cells.Item(78, 2).Value = "CONST_AIR_CO2_MIN"
$ws.Cells.Item(78, 3).Value = 1
$ws.Cells.Item(78, 4).Value = "WORD"
$ws.Cells.Item(78, 4).HorizontalAlignment = -4131
$ws.Cells.Item(78, 5).Value = 0
$ws.Cells.Item(78, 6).Value = "x"
$ws.Cells.Item(78, 7).Value = "For range checking"

$ws.Cells.Item(79, 2).Value = "CONST_WATER_RANGE_MAX"
$ws.Cells.Item(79, 3).Value = 1
$ws.Cells.Item(79, 4).Value = "WORD"
$ws.Cells.Item(79, 4).HorizontalAlignment = -4131
$ws.Cells.Item(79, 5).Value = 0
$ws.Cells.Item(79, 6).Value = "x"
$ws.Cells.Item(79, 7).Value = "Checking water temperature range"

$ws.Cells.Item(80, 2).Value = "CONST_WATER_RANGE_MIN"
$ws.Cells.Item(80, 3).Value = 1
$ws.Cells.Item(80, 4).Value = "WORD"
$ws.Cells.Item(80, 4).HorizontalAlignment = -4131
$ws.Cells.Item(80, 5).Value = 0
$ws.Cells.Item(80, 6).Value = "x"
$ws.Cells.Item(80, 7).Value = "Checking water temperature range"

# -- The pre-existing "tmpReal" row (old row 71) is now row 81; bump its
#    addr_offset from 1 to 2 to make room and give it left alignment.
$ws.Cells.Item(81, 3).Value = 2
$ws.Cells.Item(81, 4).HorizontalAlignment = -4131

# -- Insert 1 new row at 64 (global alarm-clear flag), right after the
#    last ventilation status bit (old row 62, "stat_vent_air_temp1").
$ws.Rows.Item(64).Insert()
$ws.Cells.Item(64, 2).Value = "stat_alarm_clear"
$ws.Cells.Item(64, 3).Value = 1
$ws.Cells.Item(64, 4).Value = "BOOL"
$ws.Cells.Item(64, 5).Value = $true
$ws.Cells.Item(64, 6).Value = "x"
$ws.Cells.Item(64, 7).Value = "Global alarm indicator"

# -- Insert 1 new row at 46 (new VFD "other" fault bit), right before
#    the old row 46 ("stat_water_shelf_valvefault1"), and move the
#    "Status for water system" comment that used to sit on that row to
#    this new bit instead.
$ws.Rows.Item(46).Insert()
$ws.Cells.Item(46, 2).Value = "stat_water_vfd_other1"
$ws.Cells.Item(46, 3).Value = 1
$ws.Cells.Item(46, 4).Value = "BOOL"
$ws.Cells.Item(46, 5).Value = $true
$ws.Cells.Item(46, 6).Value = "x"
$ws.Cells.Item(46, 7).Value = "Other VFD error"

# ---------------------------------------------------------------------
# 2) "IO Mapping" worksheet (sheet6.xml)
# ---------------------------------------------------------------------
$io = $wb.Worksheets.Item("IO Mapping")

# Widen column C (type) to fit "ARRAY [10] OF WORD".
$io.Columns.Item(3).ColumnWidth = 20.109375

# Insert the three new fault-record rows bottom-up so row numbers stay
# stable across inserts.

# After "MB_pump2_rVfdMode" (old row 148) -> new row 151
$io.Rows.Item(151).Insert()
$io.Cells.Item(151, 1).Value = "MB_pump2_faultRecord"
$io.Cells.Item(151, 2).Value = "D25560"
$io.Cells.Item(151, 3).Value = "ARRAY [10] OF WORD"
$io.Cells.Item(151, 3).HorizontalAlignment = -4108
$io.Cells.Item(151, 4).Value = "[10(-1)]"
$io.Cells.Item(151, 5).Value = "x"

# After "MB_pump1_rVfdMode" (old row 138) -> new row 140
$io.Rows.Item(140).Insert()
$io.Cells.Item(140, 1).Value = "MB_pump1_faultRecord"
$io.Cells.Item(140, 2).Value = "D25360"
$io.Cells.Item(140, 3).Value = "ARRAY [10] OF WORD"
$io.Cells.Item(140, 3).HorizontalAlignment = -4108
$io.Cells.Item(140, 4).Value = "[10(-1)]"
$io.Cells.Item(140, 5).Value = "x"

# After "MB_pump0_rVfdMode" (old row 128) -> new row 129
$io.Rows.Item(129).Insert()
$io.Cells.Item(129, 1).Value = "MB_pump0_faultRecord"
$io.Cells.Item(129, 2).Value = "D25160"
$io.Cells.Item(129, 3).Value = "ARRAY [10] OF WORD"
$io.Cells.Item(129, 3).HorizontalAlignment = -4108
$io.Cells.Item(129, 4).Value = "[10(-1)]"
$io.Cells.Item(129, 5).Value = "x"

# ---------------------------------------------------------------------
# 3) View / selection cosmetics
# ---------------------------------------------------------------------

# Constants: selection moves to D47
$ws.Activate()
$ws.Range("D47").Select()

# Pump sheet: selection moves to E8
$pump = $wb.Worksheets.Item("Pump")
$pump.Activate()
$pump.Range("E8").Select()

# IO Mapping: selection moves to C133
$io.Activate()
$io.Range("C133").Select()

# HMI Internal: keep selection at H10 (already set), re-activate so the
# sheet's scroll/view state is refreshed.
$hmi = $wb.Worksheets.Item("HMI Internal")
$hmi.Activate()
$hmi.Range("H10").Select()

# Re-activate Constants sheet last, matching tabSelected="1" in the
# original file.
$ws.Activate()
$ws.Range("D47").Select()

Write-Host "edit complete"
